# Report update @ 30-Aug-2025 23:59:41
# Append two new sighting rows (row 29 and row 30) to the "Sightings Data" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sightings Data")

# Copy formatting from the last existing data row (28) down into the two
# new rows so column A keeps the numbering style and column C keeps the
# date-number-format style, matching how the rest of the table is styled.
$ws.Range("A28:M28").Copy()
$ws.Range("A29:M29").PasteSpecial(-4122)
$ws.Range("A28:M28").Copy()
$ws.Range("A30:M30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 29
$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = 15511
$ws.Cells.Item(29, 3).Value = 45899
$ws.Cells.Item(29, 4).Value = "Seychelles International Airport"
$ws.Cells.Item(29, 5).Value = "Seychelles Islands Development Co"
$ws.Cells.Item(29, 6).Value = "S7DES"
$ws.Cells.Item(29, 7).Value = "SEZ"
$ws.Cells.Item(29, 8).Value = "DES"
$ws.Cells.Item(29, 9).Value = "S7-DES"
$ws.Cells.Item(29, 10).Value = "UE-284"
$ws.Cells.Item(29, 11).Value = "1900D"
$ws.Cells.Item(29, 12).Value = "Beech"
$ws.Cells.Item(29, 13).Value = "Other"

# Row 30
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = 15513
$ws.Cells.Item(30, 3).Value = 45899
$ws.Cells.Item(30, 4).Value = "Seychelles International Airport"
$ws.Cells.Item(30, 5).Value = "Seychelles Islands Development Co"
$ws.Cells.Item(30, 6).Value = "S7ALP"
$ws.Cells.Item(30, 7).Value = "DES"
$ws.Cells.Item(30, 8).Value = "SEZ"
$ws.Cells.Item(30, 9).Value = "S7-ALP"
$ws.Cells.Item(30, 10).Value = "UE-397"
$ws.Cells.Item(30, 11).Value = "1900D"
$ws.Cells.Item(30, 12).Value = "Beech"
$ws.Cells.Item(30, 13).Value = "Other"
